# Updated cryptos list on Sat Apr 29 22:20:34 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay a text value (matches the source sheet, where
    # every Price/Volume column is stored as a literal string, even when it
    # looks numeric), then restore the cell's original (default) style so we
    # don't leave a stray number-format behind.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Rows where only the Volume(1h) column (E) changes
Set-TextValue "E11" "  -0.35%  "
Set-TextValue "E17" "  +0.56%  "
Set-TextValue "E18" "  +0.97%  "
Set-TextValue "E38" "  +0.74%  "
Set-TextValue "E40" "  +1.81%  "

# Rows where Price (D) and Volume(1h) (E) both change
Set-TextValue "D2" "29.449.75"
Set-TextValue "E2" "  -0.09%  "

Set-TextValue "D3" "1.919.58"
Set-TextValue "E3" "  +0.89%  "

Set-TextValue "D5" "325.61"
Set-TextValue "E5" "  -0.07%  "

Set-TextValue "D6" "1.009"
Set-TextValue "E6" "  +0.60%  "

Set-TextValue "D7" "0.4829"
Set-TextValue "E7" "  +0.83%  "

Set-TextValue "D8" "0.4072"
Set-TextValue "E8" "  +0.16%  "

Set-TextValue "D9" "0.08237"
Set-TextValue "E9" "  +2.04%  "

Set-TextValue "D10" "1.013"
Set-TextValue "E10" "  +1.05%  "

Set-TextValue "D12" "1.925.94"
Set-TextValue "E12" "  +1.23%  "

Set-TextValue "D13" "6.071"
Set-TextValue "E13" "  +2.05%  "

Set-TextValue "D14" "7.249"
Set-TextValue "E14" "  +2.50%  "

Set-TextValue "D15" "91.76"
Set-TextValue "E15" "  +2.14%  "

Set-TextValue "D16" "0.06869"
Set-TextValue "E16" "  +2.59%  "

Set-TextValue "D19" "17.59"
Set-TextValue "E19" "  -0.13%  "

Set-TextValue "D20" "1.009"
Set-TextValue "E20" "  +0.60%  "

Set-TextValue "D21" "29.470.65"
Set-TextValue "E21" "  -0.03%  "

Set-TextValue "D22" "5.673"
Set-TextValue "E22" "  +2.40%  "

Set-TextValue "D23" "11.74"
Set-TextValue "E23" "  -0.24%  "

Set-TextValue "D24" "2.183"
Set-TextValue "E24" "  +0.73%  "

Set-TextValue "D25" "2.146.81"
Set-TextValue "E25" "  +1.19%  "

Set-TextValue "D26" "6.638"
Set-TextValue "E26" "  +9.10%  "

Set-TextValue "D27" "156.15"
Set-TextValue "E27" "  +0.94%  "

Set-TextValue "D28" "20.06"
Set-TextValue "E28" "  +1.14%  "

Set-TextValue "D29" "2.117"
Set-TextValue "E29" "  +1.09%  "

Set-TextValue "D30" "120.93"
Set-TextValue "E30" "  +2.22%  "

Set-TextValue "D31" "1.018"
Set-TextValue "E31" "  -1.80%  "

Set-TextValue "D32" "0.09622"
Set-TextValue "E32" "  +1.16%  "

Set-TextValue "D33" "5.658"
Set-TextValue "E33" "  +4.17%  "

Set-TextValue "D34" "3.552"
Set-TextValue "E34" "  +0.20%  "

Set-TextValue "D35" "1.375"
Set-TextValue "E35" "  -1.28%  "

Set-TextValue "D36" "0.02288"
Set-TextValue "E36" "  +1.66%  "

Set-TextValue "D37" "0.06107"
Set-TextValue "E37" "  +0.53%  "

Set-TextValue "D39" "8.069"
Set-TextValue "E39" "  +1.90%  "

Set-TextValue "D41" "10.86"
Set-TextValue "E41" "  +6.46%  "

Set-TextValue "D42" "0.1848"
Set-TextValue "E42" "  +0.16%  "

Set-TextValue "D43" "1.283"
Set-TextValue "E43" "  +0.26%  "

Set-TextValue "D44" "2.390"
Set-TextValue "E44" "  -1.13%  "

# Row 45 and 46 swap places (EnergySwap <-> Cronos) with new values
Set-TextValue "B45" "Cronos"
Set-TextValue "C45" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D45" "0.07607"
Set-TextValue "E45" "  -2.40%  "

Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "12.44"
Set-TextValue "E46" "  +1.65%  "

Set-TextValue "D47" "0.5604"
Set-TextValue "E47" "  +1.35%  "

Set-TextValue "D48" "1.954"
Set-TextValue "E48" "  +1.67%  "

Set-TextValue "D49" "118.70"
Set-TextValue "E49" "  +4.43%  "

Set-TextValue "D50" "2.428"
Set-TextValue "E50" "  +3.77%  "

Set-TextValue "D51" "72.35"
Set-TextValue "E51" "  +0.25%  "
